# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1455
    $ws.Range("F9").Value = 240
}
